$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to
# text format first, otherwise Excel auto-converts them to numeric cells
# (losing the original text formatting, e.g. trailing zeros).

$ws.Range('D2').Value = '66.138.63'
$ws.Range('E2').Value = '  +4.83%  '
$ws.Range('D3').Value = '3.829.66'
$ws.Range('E3').Value = '  +9.98%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '427.15'
$ws.Range('E5').Value = '  +9.69%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '130.47'
$ws.Range('E6').Value = '  +8.11%  '
$ws.Range('D7').Value = '3.828.13'
$ws.Range('E7').Value = '  +10.23%  '
$ws.Range('E8').Value = '  +4.31%  '
$ws.Range('E9').Value = '  +0.04%  '
$ws.Range('E10').Value = '  +8.17%  '
$ws.Range('E11').Value = '  +3.61%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000336'
$ws.Range('E12').Value = '  +1.94%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '41.57'
$ws.Range('E13').Value = '  +7.83%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '10.43'
$ws.Range('E14').Value = '  +14.10%  '
$ws.Range('D15').Value = '4.430.16'
$ws.Range('E15').Value = '  +9.91%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '15.57'
$ws.Range('E16').Value = '  +24.65%  '
$ws.Range('B17').Value = 'TRON'
$ws.Range('C17').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.138'
$ws.Range('E17').Value = '  +1.47%  '
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').Value = '3.820.66'
$ws.Range('E18').Value = '  +10.01%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '19.98'
$ws.Range('E19').Value = '  +7.13%  '
$ws.Range('E20').Value = '  +8.91%  '
$ws.Range('D21').Value = '66.324.46'
$ws.Range('E21').Value = '  +5.17%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '414.71'
$ws.Range('E22').Value = '  +5.47%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '15.10'
$ws.Range('E23').Value = '  +8.85%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '84.94'
$ws.Range('E24').Value = '  +5.24%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.12'
$ws.Range('E25').Value = '  +9.24%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '37.24'
$ws.Range('E26').Value = '  +12.00%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.17'
$ws.Range('E27').Value = '  +16.17%  '
$ws.Range('E28').Value = '  +10.59%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.51'
$ws.Range('E29').Value = '  +41.74%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.40'
$ws.Range('E30').Value = '  +3.81%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '13.88'
$ws.Range('E31').Value = '  +18.08%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '716.33'
$ws.Range('E32').Value = '  +6.82%  '
$ws.Range('E33').Value = '  +14.24%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.77'
$ws.Range('E34').Value = '  +6.56%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.998'
$ws.Range('E35').Value = '  -0.21%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '38.82'
$ws.Range('E36').Value = '  +6.59%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.68'
$ws.Range('E37').Value = '  +42.69%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.149'
$ws.Range('E38').Value = '  +0.55%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '55.60'
$ws.Range('E39').Value = '  +4.39%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0469'
$ws.Range('E40').Value = '  +8.40%  '
$ws.Range('D41').Value = '0.0₃0729'
$ws.Range('E41').Value = '  +14.83%  '
$ws.Range('E42').Value = '  +6.61%  '
$ws.Range('E43').Value = '  +0.28%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.29'
$ws.Range('E44').Value = '  +7.66%  '
$ws.Range('E45').Value = '  +4.57%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.38'
$ws.Range('E46').Value = '  +11.07%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.320'
$ws.Range('E47').Value = '  +17.18%  '
$ws.Range('B48').Value = 'Fetch.AI'
$ws.Range('C48').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.41'
$ws.Range('E48').Value = '  +45.86%  '
$ws.Range('B49').Value = 'WEMIXToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.63'
$ws.Range('E49').Value = '  +7.99%  '
$ws.Range('E50').Value = '  +6.22%  '
$ws.Range('E51').Value = '  +3.92%  '

# Restore default (Normal) style on the cells we force-formatted above so
# we do not leave a stray number-format override on them.
$ws.Range('D4').Style = 'Normal'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D12').Style = 'Normal'
$ws.Range('D13').Style = 'Normal'
$ws.Range('D14').Style = 'Normal'
$ws.Range('D16').Style = 'Normal'
$ws.Range('D17').Style = 'Normal'
$ws.Range('D19').Style = 'Normal'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D23').Style = 'Normal'
$ws.Range('D24').Style = 'Normal'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D27').Style = 'Normal'
$ws.Range('D29').Style = 'Normal'
$ws.Range('D30').Style = 'Normal'
$ws.Range('D31').Style = 'Normal'
$ws.Range('D32').Style = 'Normal'
$ws.Range('D34').Style = 'Normal'
$ws.Range('D35').Style = 'Normal'
$ws.Range('D36').Style = 'Normal'
$ws.Range('D37').Style = 'Normal'
$ws.Range('D38').Style = 'Normal'
$ws.Range('D39').Style = 'Normal'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D44').Style = 'Normal'
$ws.Range('D46').Style = 'Normal'
$ws.Range('D47').Style = 'Normal'
$ws.Range('D48').Style = 'Normal'
$ws.Range('D49').Style = 'Normal'
